$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.2057716666666667
$ws.Range("H2").Value2 = 0.6173149999999999
$ws.Range("I2").Value2 = 0.01089677771948535
$ws.Range("J2").Value2 = 0.01089677771948535
$ws.Range("M2").Value2 = 127.3992563333333
$ws.Range("N2").Value2 = 382.197769
$ws.Range("O2").Value2 = 0.4838549810199306
$ws.Range("P2").Value2 = 0.4838549810199307
$ws.Range("Q2").Value2 = 26.21515730780389
$ws.Range("R2").Value2 = 235.936415770235
$ws.Range("S2").Value2 = 0.005272460176639985
$ws.Range("T2").Value2 = 0.005272460176639986
$ws.Range("G3").Value2 = 0.2057716666666667
$ws.Range("H3").Value2 = 0.6173149999999999
$ws.Range("I3").Value2 = 0.01089677771948535
$ws.Range("J3").Value2 = 0.01089677771948535
$ws.Range("M3").Value2 = 59.36586533333332
$ws.Range("O3").Value2 = 0.2254681108101269
$ws.Range("P3").Value2 = 0.2254681108101269
$ws.Range("Q3").Value2 = 12.21581305274889
$ws.Range("R3").Value2 = 109.94231747474
$ws.Range("S3").Value2 = 0.002456875886330243
$ws.Range("T3").Value2 = 0.002456875886330244
$ws.Range("G4").Value2 = 0.2057716666666667
$ws.Range("H4").Value2 = 0.6173149999999999
$ws.Range("I4").Value2 = 0.01089677771948535
$ws.Range("J4").Value2 = 0.01089677771948535
$ws.Range("M4").Value2 = 16.63275166666667
$ws.Range("N4").Value2 = 49.898255
$ws.Range("O4").Value2 = 0.06317022542837675
$ws.Range("P4").Value2 = 0.06317022542837675
$ws.Range("Q4").Value2 = 3.422549031702778
$ws.Range("R4").Value2 = 30.802941285325
$ws.Range("S4").Value2 = 0.0006883519049828024
$ws.Range("T4").Value2 = 0.0006883519049828024
$ws.Range("G5").Value2 = 0.2057716666666667
$ws.Range("H5").Value2 = 0.6173149999999999
$ws.Range("I5").Value2 = 0.01089677771948535
$ws.Range("J5").Value2 = 0.01089677771948535
$ws.Range("M5").Value2 = 59.90262233333334
$ws.Range("N5").Value2 = 179.707867
$ws.Range("O5").Value2 = 0.2275066827415657
$ws.Range("P5").Value2 = 0.2275066827415658
$ws.Range("Q5").Value2 = 12.32626243523389
$ws.Range("R5").Value2 = 110.936361917105
$ws.Range("S5").Value2 = 0.002479089751532315
$ws.Range("T5").Value2 = 0.002479089751532315
$ws.Range("I6").Value2 = 0.01769706320706529
$ws.Range("J6").Value2 = 0.01769706320706529
$ws.Range("M6").Value2 = 127.3992563333333
$ws.Range("N6").Value2 = 382.197769
$ws.Range("O6").Value2 = 0.4838549810199306
$ws.Range("P6").Value2 = 0.4838549810199307
$ws.Range("Q6").Value2 = 42.57509034343011
$ws.Range("R6").Value2 = 383.175813090871
$ws.Range("S6").Value2 = 0.008562812182163088
$ws.Range("T6").Value2 = 0.008562812182163091
$ws.Range("I7").Value2 = 0.01769706320706529
$ws.Range("J7").Value2 = 0.01769706320706529
$ws.Range("M7").Value2 = 59.36586533333332
$ws.Range("O7").Value2 = 0.2254681108101269
$ws.Range("P7").Value2 = 0.2254681108101269
$ws.Range("S7").Value2 = 0.003990123408184415
$ws.Range("T7").Value2 = 0.003990123408184417
$ws.Range("I8").Value2 = 0.01769706320706529
$ws.Range("J8").Value2 = 0.01769706320706529
$ws.Range("M8").Value2 = 16.63275166666667
$ws.Range("N8").Value2 = 49.898255
$ws.Range("O8").Value2 = 0.06317022542837675
$ws.Range("P8").Value2 = 0.06317022542837675
$ws.Range("Q8").Value2 = 5.558438292727222
$ws.Range("R8").Value2 = 50.025944634545
$ws.Range("S8").Value2 = 0.001117927472210546
$ws.Range("T8").Value2 = 0.001117927472210546
$ws.Range("I9").Value2 = 0.01769706320706529
$ws.Range("J9").Value2 = 0.01769706320706529
$ws.Range("M9").Value2 = 59.90262233333334
$ws.Range("N9").Value2 = 179.707867
$ws.Range("O9").Value2 = 0.2275066827415657
$ws.Range("P9").Value2 = 0.2275066827415658
$ws.Range("Q9").Value2 = 20.01863771462811
$ws.Range("R9").Value2 = 180.167739431653
$ws.Range("S9").Value2 = 0.004026200144507238
$ws.Range("T9").Value2 = 0.00402620014450724
$ws.Range("G10").Value2 = 0.4895776666666666
$ws.Range("H10").Value2 = 1.468733
$ws.Range("I10").Value2 = 0.02592591631545138
$ws.Range("J10").Value2 = 0.02592591631545138
$ws.Range("M10").Value2 = 127.3992563333333
$ws.Range("N10").Value2 = 382.197769
$ws.Range("O10").Value2 = 0.4838549810199306
$ws.Range("P10").Value2 = 0.4838549810199307
$ws.Range("Q10").Value2 = 62.37183065074188
$ws.Range("R10").Value2 = 561.3464758566769
$ws.Range("S10").Value2 = 0.01254438374673704
$ws.Range("T10").Value2 = 0.01254438374673704
$ws.Range("G11").Value2 = 0.4895776666666666
$ws.Range("H11").Value2 = 1.468733
$ws.Range("I11").Value2 = 0.02592591631545138
$ws.Range("J11").Value2 = 0.02592591631545138
$ws.Range("M11").Value2 = 59.36586533333332
$ws.Range("O11").Value2 = 0.2254681108101269
$ws.Range("P11").Value2 = 0.2254681108101269
$ws.Range("Q11").Value2 = 29.06420182954088
$ws.Range("R11").Value2 = 261.577816465868
$ws.Range("S11").Value2 = 0.005845467372666268
$ws.Range("T11").Value2 = 0.005845467372666269
$ws.Range("G12").Value2 = 0.4895776666666666
$ws.Range("H12").Value2 = 1.468733
$ws.Range("I12").Value2 = 0.02592591631545138
$ws.Range("J12").Value2 = 0.02592591631545138
$ws.Range("M12").Value2 = 16.63275166666667
$ws.Range("N12").Value2 = 49.898255
$ws.Range("O12").Value2 = 0.06317022542837675
$ws.Range("P12").Value2 = 0.06317022542837675
$ws.Range("Q12").Value2 = 8.143023751212777
$ws.Range("R12").Value2 = 73.28721376091499
$ws.Range("S12").Value2 = 0.001637745978084294
$ws.Range("T12").Value2 = 0.001637745978084294
$ws.Range("G13").Value2 = 0.4895776666666666
$ws.Range("H13").Value2 = 1.468733
$ws.Range("I13").Value2 = 0.02592591631545138
$ws.Range("J13").Value2 = 0.02592591631545138
$ws.Range("M13").Value2 = 59.90262233333334
$ws.Range("N13").Value2 = 179.707867
$ws.Range("O13").Value2 = 0.2275066827415657
$ws.Range("P13").Value2 = 0.2275066827415658
$ws.Range("Q13").Value2 = 29.32698606916789
$ws.Range("R13").Value2 = 263.942874622511
$ws.Range("S13").Value2 = 0.005898319217963781
$ws.Range("T13").Value2 = 0.005898319217963782
$ws.Range("G14").Value2 = 17.85418133333333
$ws.Range("H14").Value2 = 53.562544
$ws.Range("I14").Value2 = 0.9454802427579979
$ws.Range("J14").Value2 = 0.945480242757998
$ws.Range("M14").Value2 = 127.3992563333333
$ws.Range("N14").Value2 = 382.197769
$ws.Range("O14").Value2 = 0.4838549810199306
$ws.Range("P14").Value2 = 0.4838549810199307
$ws.Range("Q14").Value2 = 2274.609424307148
$ws.Range("R14").Value2 = 20471.48481876434
$ws.Range("S14").Value2 = 0.4574753249143905
$ws.Range("T14").Value2 = 0.4574753249143906
$ws.Range("G15").Value2 = 17.85418133333333
$ws.Range("H15").Value2 = 53.562544
$ws.Range("I15").Value2 = 0.9454802427579979
$ws.Range("J15").Value2 = 0.945480242757998
$ws.Range("M15").Value2 = 59.36586533333332
$ws.Range("O15").Value2 = 0.2254681108101269
$ws.Range("P15").Value2 = 0.2254681108101269
$ws.Range("Q15").Value2 = 1059.92892467158
$ws.Range("R15").Value2 = 9539.360322044224
$ws.Range("S15").Value2 = 0.2131756441429459
$ws.Range("T15").Value2 = 0.213175644142946
$ws.Range("G16").Value2 = 17.85418133333333
$ws.Range("H16").Value2 = 53.562544
$ws.Range("I16").Value2 = 0.9454802427579979
$ws.Range("J16").Value2 = 0.945480242757998
$ws.Range("M16").Value2 = 16.63275166666667
$ws.Range("N16").Value2 = 49.898255
$ws.Range("O16").Value2 = 0.06317022542837675
$ws.Range("P16").Value2 = 0.06317022542837675
$ws.Range("Q16").Value2 = 296.9641643289689
$ws.Range("R16").Value2 = 2672.67747896072
$ws.Range("S16").Value2 = 0.0597262000730991
$ws.Range("T16").Value2 = 0.0597262000730991
$ws.Range("G17").Value2 = 17.85418133333333
$ws.Range("H17").Value2 = 53.562544
$ws.Range("I17").Value2 = 0.9454802427579979
$ws.Range("J17").Value2 = 0.945480242757998
$ws.Range("M17").Value2 = 59.90262233333334
$ws.Range("N17").Value2 = 179.707867
$ws.Range("O17").Value2 = 0.2275066827415657
$ws.Range("P17").Value2 = 0.2275066827415658
$ws.Range("Q17").Value2 = 1069.512281481517
$ws.Range("R17").Value2 = 9625.610533333649
$ws.Range("S17").Value2 = 0.2151030736275624
$ws.Range("T17").Value2 = 0.2151030736275625
